$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("France Ligue 1")

# Row 122
$ws.Cells.Item(122, 2).Value = 6847939
$ws.Cells.Item(122, 5).Value = "Nantes"
$ws.Cells.Item(122, 6).Value = "Le Havre"
$ws.Cells.Item(122, 7).Value = 0
$ws.Cells.Item(122, 8).Value = 0
$ws.Cells.Item(122, 9).Value = "D"
$ws.Cells.Item(122, 11).Value = 3.3
$ws.Cells.Item(122, 12).Value = 3.4
$ws.Cells.Item(122, 13).Value = 1.95
$ws.Cells.Item(122, 14).Value = 3.25
$ws.Cells.Item(122, 15).Value = 4.333
$ws.Cells.Item(122, 16).Value = -0.5
$ws.Cells.Item(122, 19).Value = 2
$ws.Cells.Item(122, 20).Value = 1.825
$ws.Cells.Item(122, 21).Value = 2.025
$ws.Cells.Item(122, 23).Value = 2.25
$ws.Cells.Item(122, 24).Value = -1
$ws.Cells.Item(122, 27).Value = -1
$ws.Cells.Item(122, 28).Value = 1.025

# Row 123
$ws.Cells.Item(123, 2).Value = 6847940
$ws.Cells.Item(123, 5).Value = "Montpellier"
$ws.Cells.Item(123, 6).Value = "Brest"
$ws.Cells.Item(123, 7).Value = 1
$ws.Cells.Item(123, 8).Value = 3
$ws.Cells.Item(123, 9).Value = "A"
$ws.Cells.Item(123, 11).Value = 3.4
$ws.Cells.Item(123, 12).Value = 3.25
$ws.Cells.Item(123, 13).Value = 2.25
$ws.Cells.Item(123, 14).Value = 3.3
$ws.Cells.Item(123, 15).Value = 3.25
$ws.Cells.Item(123, 16).Value = -0.25
$ws.Cells.Item(123, 19).Value = 2.5
$ws.Cells.Item(123, 20).Value = 2.05
$ws.Cells.Item(123, 21).Value = 1.8
$ws.Cells.Item(123, 23).Value = -1
$ws.Cells.Item(123, 24).Value = 2.25
$ws.Cells.Item(123, 27).Value = 1.05
$ws.Cells.Item(123, 28).Value = -1

# Row 143
$ws.Cells.Item(143, 2).Value = 6847951
$ws.Cells.Item(143, 5).Value = "Strasbourg"
$ws.Cells.Item(143, 6).Value = "Le Havre"
$ws.Cells.Item(143, 7).Value = 2
$ws.Cells.Item(143, 8).Value = 1
$ws.Cells.Item(143, 9).Value = "H"
$ws.Cells.Item(143, 10).Value = 2.1
$ws.Cells.Item(143, 12).Value = 3.75
$ws.Cells.Item(143, 13).Value = 2.375
$ws.Cells.Item(143, 14).Value = 3.1
$ws.Cells.Item(143, 15).Value = 3.4
$ws.Cells.Item(143, 16).Value = -0.25
$ws.Cells.Item(143, 17).Value = 2.02
$ws.Cells.Item(143, 18).Value = 1.88
$ws.Cells.Item(143, 19).Value = 2
$ws.Cells.Item(143, 20).Value = 1.86
$ws.Cells.Item(143, 21).Value = 2.04
$ws.Cells.Item(143, 22).Value = 1.375
$ws.Cells.Item(143, 23).Value = -1
$ws.Cells.Item(143, 25).Value = 1.02
$ws.Cells.Item(143, 27).Value = 0.8600000000000001
$ws.Cells.Item(143, 28).Value = -1

# Row 144
$ws.Cells.Item(144, 2).Value = 6847956
$ws.Cells.Item(144, 5).Value = "Clermont Foot"
$ws.Cells.Item(144, 6).Value = "Lille"
$ws.Cells.Item(144, 7).Value = 0
$ws.Cells.Item(144, 8).Value = 0
$ws.Cells.Item(144, 9).Value = "D"
$ws.Cells.Item(144, 10).Value = 4
$ws.Cells.Item(144, 12).Value = 2.05
$ws.Cells.Item(144, 13).Value = 4.333
$ws.Cells.Item(144, 14).Value = 3.5
$ws.Cells.Item(144, 15).Value = 1.85
$ws.Cells.Item(144, 16).Value = 0.5
$ws.Cells.Item(144, 17).Value = 1.975
$ws.Cells.Item(144, 18).Value = 1.875
$ws.Cells.Item(144, 19).Value = 2.5
$ws.Cells.Item(144, 20).Value = 1.95
$ws.Cells.Item(144, 21).Value = 1.9
$ws.Cells.Item(144, 22).Value = -1
$ws.Cells.Item(144, 23).Value = 2.5
$ws.Cells.Item(144, 25).Value = 0.9750000000000001
$ws.Cells.Item(144, 27).Value = -1
$ws.Cells.Item(144, 28).Value = 0.8999999999999999

# Row 178
$ws.Cells.Item(178, 2).Value = 6847992
$ws.Cells.Item(178, 5).Value = "Clermont Foot"
$ws.Cells.Item(178, 6).Value = "Strasbourg"
$ws.Cells.Item(178, 7).Value = 1
$ws.Cells.Item(178, 8).Value = 1
$ws.Cells.Item(178, 10).Value = 2.8
$ws.Cells.Item(178, 11).Value = 3.6
$ws.Cells.Item(178, 12).Value = 2.3
$ws.Cells.Item(178, 13).Value = 3
$ws.Cells.Item(178, 14).Value = 3.3
$ws.Cells.Item(178, 15).Value = 2.375
$ws.Cells.Item(178, 17).Value = 1.8
$ws.Cells.Item(178, 18).Value = 2.05
$ws.Cells.Item(178, 19).Value = 2.25
$ws.Cells.Item(178, 20).Value = 1.85
$ws.Cells.Item(178, 21).Value = 2
$ws.Cells.Item(178, 23).Value = 2.3
$ws.Cells.Item(178, 25).Value = 0.4
$ws.Cells.Item(178, 27).Value = -0.5
$ws.Cells.Item(178, 28).Value = 0.5

# Row 179
$ws.Cells.Item(179, 2).Value = 6847989
$ws.Cells.Item(179, 5).Value = "Lorient"
$ws.Cells.Item(179, 6).Value = "Le Havre"
$ws.Cells.Item(179, 7).Value = 3
$ws.Cells.Item(179, 8).Value = 3
$ws.Cells.Item(179, 10).Value = 2.5
$ws.Cells.Item(179, 11).Value = 3.25
$ws.Cells.Item(179, 12).Value = 2.8
$ws.Cells.Item(179, 13).Value = 3.1
$ws.Cells.Item(179, 14).Value = 3.1
$ws.Cells.Item(179, 15).Value = 2.5
$ws.Cells.Item(179, 17).Value = 1.81
$ws.Cells.Item(179, 18).Value = 2.125
$ws.Cells.Item(179, 19).Value = 2
$ws.Cells.Item(179, 20).Value = 2
$ws.Cells.Item(179, 21).Value = 1.9
$ws.Cells.Item(179, 23).Value = 2.1
$ws.Cells.Item(179, 25).Value = 0.405
$ws.Cells.Item(179, 27).Value = 1
$ws.Cells.Item(179, 28).Value = -1

# Row 231
$ws.Cells.Item(231, 2).Value = 6848047
$ws.Cells.Item(231, 5).Value = "Le Havre"
$ws.Cells.Item(231, 6).Value = "Toulouse"
$ws.Cells.Item(231, 12).Value = 2.9
$ws.Cells.Item(231, 13).Value = 3
$ws.Cells.Item(231, 14).Value = 3
$ws.Cells.Item(231, 15).Value = 2.6
$ws.Cells.Item(231, 16).Value = 0
$ws.Cells.Item(231, 17).Value = 2.05
$ws.Cells.Item(231, 18).Value = 1.8
$ws.Cells.Item(231, 20).Value = 2.1
$ws.Cells.Item(231, 21).Value = 1.775
$ws.Cells.Item(231, 22).Value = 2
$ws.Cells.Item(231, 25).Value = 1.05
$ws.Cells.Item(231, 28).Value = 0.7749999999999999

# Row 232
$ws.Cells.Item(232, 2).Value = 6848041
$ws.Cells.Item(232, 5).Value = "Strasbourg"
$ws.Cells.Item(232, 6).Value = "Monaco"
$ws.Cells.Item(232, 7).Value = 0
$ws.Cells.Item(232, 8).Value = 1
$ws.Cells.Item(232, 9).Value = "A"
$ws.Cells.Item(232, 10).Value = 3.4
$ws.Cells.Item(232, 11).Value = 3.6
$ws.Cells.Item(232, 12).Value = 2.05
$ws.Cells.Item(232, 13).Value = 3.3
$ws.Cells.Item(232, 14).Value = 3.75
$ws.Cells.Item(232, 15).Value = 2.05
$ws.Cells.Item(232, 16).Value = 0.5
$ws.Cells.Item(232, 17).Value = 1.84
$ws.Cells.Item(232, 18).Value = 2.09
$ws.Cells.Item(232, 19).Value = 3
$ws.Cells.Item(232, 20).Value = 2
$ws.Cells.Item(232, 21).Value = 1.9
$ws.Cells.Item(232, 22).Value = -1
$ws.Cells.Item(232, 24).Value = 1.05
$ws.Cells.Item(232, 25).Value = -1
$ws.Cells.Item(232, 26).Value = 1.09
$ws.Cells.Item(232, 28).Value = 0.8999999999999999

# Row 233
$ws.Cells.Item(233, 2).Value = 6848048
$ws.Cells.Item(233, 5).Value = "Metz"
$ws.Cells.Item(233, 6).Value = "Clermont Foot"
$ws.Cells.Item(233, 7).Value = 1
$ws.Cells.Item(233, 8).Value = 0
$ws.Cells.Item(233, 9).Value = "H"
$ws.Cells.Item(233, 10).Value = 2.45
$ws.Cells.Item(233, 11).Value = 3.2
$ws.Cells.Item(233, 12).Value = 3
$ws.Cells.Item(233, 13).Value = 2.15
$ws.Cells.Item(233, 14).Value = 3.3
$ws.Cells.Item(233, 15).Value = 3.5
$ws.Cells.Item(233, 16).Value = -0.25
$ws.Cells.Item(233, 17).Value = 1.85
$ws.Cells.Item(233, 18).Value = 2
$ws.Cells.Item(233, 19).Value = 2.25
$ws.Cells.Item(233, 20).Value = 1.975
$ws.Cells.Item(233, 21).Value = 1.875
$ws.Cells.Item(233, 22).Value = 1.15
$ws.Cells.Item(233, 24).Value = -1
$ws.Cells.Item(233, 25).Value = 0.8500000000000001
$ws.Cells.Item(233, 26).Value = -1
$ws.Cells.Item(233, 28).Value = 0.875

# Row 279
$ws.Cells.Item(279, 2).Value = 7998275
$ws.Cells.Item(279, 5).Value = "Marseille"
$ws.Cells.Item(279, 6).Value = "Nice"
$ws.Cells.Item(279, 7).Value = 2
$ws.Cells.Item(279, 8).Value = 2
$ws.Cells.Item(279, 9).Value = "D"
$ws.Cells.Item(279, 10).Value = 2
$ws.Cells.Item(279, 12).Value = 3.75
$ws.Cells.Item(279, 13).Value = 2.2
$ws.Cells.Item(279, 14).Value = 3.3
$ws.Cells.Item(279, 15).Value = 3.5
$ws.Cells.Item(279, 17).Value = 1.875
$ws.Cells.Item(279, 18).Value = 1.975
$ws.Cells.Item(279, 19).Value = 2.25
$ws.Cells.Item(279, 20).Value = 2
$ws.Cells.Item(279, 21).Value = 1.85
$ws.Cells.Item(279, 22).Value = -1
$ws.Cells.Item(279, 23).Value = 2.3
$ws.Cells.Item(279, 25).Value = -0.5
$ws.Cells.Item(279, 26).Value = 0.4875
$ws.Cells.Item(279, 27).Value = 1
$ws.Cells.Item(279, 28).Value = -1

# Row 280
$ws.Cells.Item(280, 2).Value = 7998281
$ws.Cells.Item(280, 5).Value = "Monaco"
$ws.Cells.Item(280, 6).Value = "Lille"
$ws.Cells.Item(280, 7).Value = 1
$ws.Cells.Item(280, 8).Value = 0
$ws.Cells.Item(280, 9).Value = "H"
$ws.Cells.Item(280, 10).Value = 2.3
$ws.Cells.Item(280, 12).Value = 3
$ws.Cells.Item(280, 13).Value = 2.3
$ws.Cells.Item(280, 14).Value = 3.75
$ws.Cells.Item(280, 15).Value = 2.8
$ws.Cells.Item(280, 17).Value = 2.05
$ws.Cells.Item(280, 18).Value = 1.8
$ws.Cells.Item(280, 19).Value = 3
$ws.Cells.Item(280, 20).Value = 2.05
$ws.Cells.Item(280, 21).Value = 1.8
$ws.Cells.Item(280, 22).Value = 1.3
$ws.Cells.Item(280, 23).Value = -1
$ws.Cells.Item(280, 25).Value = 1.05
$ws.Cells.Item(280, 26).Value = -1
$ws.Cells.Item(280, 27).Value = -1
$ws.Cells.Item(280, 28).Value = 0.8
# Row 290
$ws.Cells.Item(290, 17).Value = 1.9
$ws.Cells.Item(290, 18).Value = 2
$ws.Cells.Item(290, 20).Value = 1.9
$ws.Cells.Item(290, 21).Value = 2

# Row 291
$ws.Cells.Item(291, 13).Value = 1.444
$ws.Cells.Item(291, 15).Value = 6.5
$ws.Cells.Item(291, 17).Value = 2
$ws.Cells.Item(291, 18).Value = 1.9
$ws.Cells.Item(291, 20).Value = 1.93
$ws.Cells.Item(291, 21).Value = 1.97

# Row 292
$ws.Cells.Item(292, 17).Value = 2.06
$ws.Cells.Item(292, 18).Value = 1.84

# Row 293
$ws.Cells.Item(293, 16).Value = -1.5
$ws.Cells.Item(293, 17).Value = 1.85
$ws.Cells.Item(293, 18).Value = 2.05
$ws.Cells.Item(293, 20).Value = 2.02
$ws.Cells.Item(293, 21).Value = 1.88

# Row 294
$ws.Cells.Item(294, 20).Value = 1.89
$ws.Cells.Item(294, 21).Value = 2.01

# Row 295
$ws.Cells.Item(295, 20).Value = 1.9
$ws.Cells.Item(295, 21).Value = 2

